$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix header typo: "Requsition ID" -> "Requisition ID"
$ws.Cells.Item(1, 3).Value = "Requisition ID"

# 2. Re-format the header row and the existing data block as Text,
#    and re-enter the numeric-looking IDs as text values so they match the
#    "t=s" (shared string) cell type used in the target workbook.
$ws.Range("A1:D1").NumberFormat = "@"
$ws.Range("A2:D6").NumberFormat = "@"

$ws.Cells.Item(2, 3).Value = "68156"
$ws.Cells.Item(4, 3).Value = "441054"

# 3. Add a new row of data: Intel / Co-op Software Intern / 123456 / Mintu
$ws.Range("A7:D7").HorizontalAlignment = -4108
$ws.Range("A7:D7").ReadingOrder = 1

$ws.Cells.Item(7, 1).Value = "Intel"
$ws.Cells.Item(7, 2).Value = "Co-op Software Intern"
$ws.Cells.Item(7, 3).Value = 123456
$ws.Cells.Item(7, 4).Value = "Mintu"
